$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("A2").Value = "2022 June 16"
    $ws.Range("D4").Value = "Name"
}

$wsAll = $wb.Worksheets.Item("All")
$wsAll.Range("D5").Value = "JS220_ep_pcb_revC"

$wsTop = $wb.Worksheets.Item("Top")
$wsTop.Range("D5").Value = "JS220_ep_pcb_revC"
